$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")
$ws.Range("C3").Value = "Universitat Jaume I"
$ws.Range("C3").Select()
